# Weekly data refresh: insert a new observation row at row 48, pushing the
# existing rows 48-154 down to 49-155.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(48).Insert()

$newDate = Get-Date -Year 2023 -Month 5 -Day 24 -Hour 0 -Minute 0 -Second 0

$ws.Range("A48").Value = 4
$ws.Range("B48").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C48").Value = "Los Lagos"
$ws.Range("D48").Value = $newDate
$ws.Range("E48").Value = 10
$ws.Range("F48").Value = "Fruta"
$ws.Range("G48").Value = 100104
$ws.Range("H48").Value = "Frutos de pepita"
$ws.Range("I48").Value = 100104003
$ws.Range("J48").Value = "Membrillo"
$ws.Range("K48").Value = "Champion"
$ws.Range("L48").Value = "Primera"
$ws.Range("M48").Value = 40
$ws.Range("N48").Value = 13000
$ws.Range("O48").Value = 14000
$ws.Range("P48").Value = 13500
$ws.Range("Q48").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R48").Value = "Región de O'Higgins"
$ws.Range("S48").Value = 750
$ws.Range("T48").Value = 18
